# Align opti with new reference: update computed values for rows 2-9 and
# remove the now-unused last data row (row 10, sample S9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns C, D, E, F, G, H)
$data = @(
    @{ Row = 2;  C = 99.27232944997668; D = 4.437436252154237; E = $true;  F = 1;    G = 0.65; H = 2 },
    @{ Row = 3;  C = 99.31815534589504; D = 4.455650981524249; E = $false; F = 0.95; G = 0.68; H = 1 },
    @{ Row = 4;  C = 98.68006956679432; D = 5.908935444679265; E = $false; F = 0.95; G = 0.68; H = 1 },
    @{ Row = 5;  C = 99.29166895405103; D = 4.058438737581729; E = $false; F = 0.95; G = 0.68; H = 1 },
    @{ Row = 6;  C = 98.82383710792956; D = 4.85395972154425;  E = $false; F = 0.95; G = 0.68; H = 1 },
    @{ Row = 7;  C = 98.79549216381972; D = 4.074358541525996; E = $false; F = 0.95; G = 0.68; H = 1 },
    @{ Row = 8;  C = 99.09193504566448; D = 6.721120573423139; E = $false; F = 0.95; G = 0.68; H = 1 },
    @{ Row = 9;  C = 98.56170310665856; D = 6.21022043795814;  E = $false; F = 0.95; G = 0.68; H = 1 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
}

# Remove the last data row (row 10, sample S9) entirely.
$ws.Rows.Item(10).Delete()
